$d = $word.ActiveDocument

foreach ($n in 1..3) {
    $old = "<id>p068r_$n</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $old, 2)
}
